$d = $word.ActiveDocument

# 1) Paragraph: "В дополнении к этому..." - replace the trailing English sentence
$old1 = "Cryptocurrencies have indeed proposed solutions to many of these challenges, as it is easy to find many projects today that offer instant transactions and virtually no fees. But, this doesn" + [char]8217 + "t change the fact that using this new technology in its current form is like trying to send an email in the 80s " + [char]8211 + " too time-consuming to setup and operate to be practical for most everyday applications."
$new1 = ". В настоящее время, криптовалюты предлагают множество способов решить эти проблемы. Сейчас можно найти довольно много проектов, которые предлагают мгновенные транзакции и практически нулевую комиссию. Но это не отменяет того факта, что использование этих технологий напоминает попытку отправить электронное письмо в 80-х годах " + [char]8212 + " трудоёмкая работа по настройке, чтобы быть практичным для большинства приложений, и осторожная эксплуатация."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2) Paragraph: "Using their new feature called the SmartCard..." full replacement
$old2 = "Using their new feature called the SmartCard, the developers of SmartCash are looking to improve upon and/or eliminate entirely the present limitations of payment cards. SmartCard is the crypto alternative to the debit card which uses the SmartCash blockchain to easily initiate and verify transactions of any size and frequency. As with many other cryptocurrencies, the fees are right around 0%; it costs just fractions of a cent to process a transaction, regardless of the amount. The payments will be instantly confirmed at the time of purchase, meaning there is no delay between authorization of the payment and the transfer of the actual funds."
$new2 = "Используя новую технологию SmartCard, разработчики SmartCash стремятся полностью устранить существующие ограничения для платёжных карт. SmartCard " + [char]8212 + " это криптовалютная альтернатива дебетовым картам, которая использует блокчейн SmartCash, чтобы легко осуществлять и проверять транзакции любого размера и с любой частотой. Как и во многих других криптовалютах, комиссия составит около 0%; вы платите лишь доли цента за обработку транзакции, независимо от вашей суммы. Платежи будут мгновенно подтверждены на момент покупки, не будет никакой задержки между авторизацией платежа и переводом средств."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3) Paragraph: "In a way, the purchases made with SmartCard..." full replacement
$old3 = "In a way, the purchases made with SmartCard combine the best of cash, card, and crypto payments. The funds are moved between parties instantly, like cash, but the transaction is authorized by the buyer with a simple code that can be either digitally stored on a smartphone or printed and carried separately as a physical QR code. This makes the payment process just as convenient as regular card payments, and far easier than having to send cryptocurrency manually from a mobile app."
$new3 = "Если сравнить, покупки через SmartCard сочетают в себе преимущества наличных, карточных и криптовалютных платежных операций. Средства мгновенно перемещаются между сторонами, как в случае с наличными деньгами, но сама транзакция требует подтверждения покупателем с помощью простого кода, который может быть либо сохранен в цифровом виде на смартфоне, либо распечатан отдельно как физический QR-код. Это делает процесс оплаты столь же удобным, как и обычные платежи по карте, но гораздо проще, чем отправка криптовалюты вручную из мобильного приложения."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# 4) Heading: "HOW IT ALL WORKS" -> "Как это работает"
$d.Content.Find.Execute("HOW IT ALL WORKS", $true, $false, $false, $false, $false, $true, 1, $false, "Как это работает", 2)

# 5) "There are just two main applications..." full replacement
$old5 = "There are just two main applications: the card app for the consumer and the Point of Sale (PoS) app for the merchant."
$new5 = "Существует два основных приложения: приложение карт для потребителей и Point of Sale (PoS) приложение для продавцов."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)

# 6) "The card app" -> "Приложение"
$d.Content.Find.Execute("The card app", $true, $false, $false, $false, $false, $true, 1, $false, "Приложение", 2)

# 7) " allows anyone to create a SmartCard..." (leading char is a non-breaking space) full replacement
$old7 = [char]160 + "allows anyone to create a SmartCard with just a few taps. Each card has an associated public address which can be used to load more funds onto the card from any wallet or exchange that offers SmartCash. Also included into each one is a QR code which can be scanned by the PoS application from either a phone screen or from a simple piece of paper. Once loaded through this address, the card works very similar to a pre-paid debit card."
$new7 = " карт позволяет создать SmartCard за несколько простых шагов. Каждая карта имеет соответствующий публичный адрес, который можно использовать для пополнения карты из любого кошелька или биржи, которая имеет в своём списке SmartCash. Кроме того, у каждой карты есть соответствующий адресу QR-код, который может быть отсканирован с экрана телефона, либо (если он распечатан) с листа бумаги. После пополнения этого адреса, SmartCard будет функционировать подобно обычной дебетовой карте."
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)

# 8) "The PoS app" -> "Приложение"
$d.Content.Find.Execute("The PoS app", $true, $false, $false, $false, $false, $true, 1, $false, "Приложение", 2)

# 9) " is equally responsible for this convenience...in the equivalent amount of" -> Russian (partial, English tail remains)
$old9 = [char]160 + "is equally responsible for this convenience. The merchant only has to specify the transaction amount in the local currency, and the app automatically initiates a transaction in the equivalent amount of"
$new9 = " PoS в такой же степени удобно. Продавцу необходимо указать лишь сумму транзакции в местной валюте, а приложение автоматически конвертирует эту сумму в эквивалентную сумму"
$d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
